# Refresh the cryptos list: update Price / Volume(1h) figures for the
# existing rows, and insert a new "OKB" row at position 39 -- shifting
# the remaining coins down by one row (InjectiveProtocol, previously the
# last row, drops off the bottom of the fixed A1:E51 range).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) store plain text in this sheet, even
# when the text looks numeric (e.g. "4.70", "0.329"). Force the Text
# format on the data range before writing so Excel does not silently
# reinterpret those values as numbers (which would drop significant
# trailing zeros and change the cell type).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '54.242.22'
$ws.Range("E2").Value = '  +1.20%  '
$ws.Range("D3").Value = '2.268.99'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '495.59'
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("D6").Value = '128.02'
$ws.Range("E6").Value = '  +2.28%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("E9").Value = '  +4.28%  '
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = '0.329'
$ws.Range("E11").Value = '  +3.37%  '
$ws.Range("D12").Value = '4.70'
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = '2.674.45'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '22.19'
$ws.Range("E14").Value = '  +4.46%  '
$ws.Range("D15").Value = '54.180.33'
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("D17").Value = '2.278.11'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '10.10'
$ws.Range("E18").Value = '  +4.58%  '
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("D20").Value = '302.71'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("E21").Value = '  +4.77%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = '61.89'
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").Value = '2.378.70'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("E26").Value = '  +1.82%  '
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("D28").Value = '171.18'
$ws.Range("E28").Value = '  +4.35%  '
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("D30").Value = '5.89'
$ws.Range("E30").Value = '  +1.70%  '
$ws.Range("D31").Value = '0.0₃0681'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '0.899'
$ws.Range("E36").Value = '  +7.33%  '
$ws.Range("E38").Value = '  +3.77%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '35.67'
$ws.Range("E39").Value = '  +1.32%  '
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '0.372'
$ws.Range("E40").Value = '  +1.02%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.41'
$ws.Range("E41").Value = '  +2.35%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '3.39'
$ws.Range("E42").Value = '  +2.82%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '126.33'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '4.77'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Value = '0.0899'
$ws.Range("E45").Value = '  +2.10%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '0.0487'
$ws.Range("E46").Value = '  +3.56%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.546'
$ws.Range("E47").Value = '  +2.05%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '237.61'
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("B49").Value = 'Polygon'
$ws.Range("C49").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D49").Value = '0.371'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0205'
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '10.76'
$ws.Range("E51").Value = '  +0.95%  '
